# Apply the "Update countries & provincias Spain" data refresh to the Pais sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Footer: refresh the "last updated" timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 26 de Agosto de 2020 a las 12:57"

# Row 4: Estados Unidos - refreshed case counts
$ws.Range("B4").Value = 5956160
$ws.Range("C4").Value = 432
$ws.Range("E4").Value = 2519000
$ws.Range("G4").Value = 17
$ws.Range("H4").Value = 182421

# Row 14: Iran - refreshed case counts
$ws.Range("B14").Value = 365606
$ws.Range("C14").Value = 2243
$ws.Range("D14").Value = 314870
$ws.Range("E14").Value = 29716
$ws.Range("G14").Value = 119
$ws.Range("H14").Value = 21020

# Row 18: Banglades - refreshed case counts
$ws.Range("B18").Value = 302147
$ws.Range("C18").Value = 2519
$ws.Range("D18").Value = 190183
$ws.Range("E18").Value = 107882
$ws.Range("G18").Value = 54
$ws.Range("H18").Value = 4082

# Row 41: country label changes from Kuwait to Rumania (re-sort), with refreshed case counts
$ws.Range("A41").Value = "Rumania"
$ws.Range("B41").Value = 81646
$ws.Range("C41").Value = 1256
$ws.Range("D41").Value = 36286
$ws.Range("E41").Value = 41939
$ws.Range("G41").Value = 54
$ws.Range("H41").Value = 3421

# Row 42: country label changes from Rumania to Kuwait (re-sort), with refreshed case counts
$ws.Range("A42").Value = "Kuwait"
$ws.Range("B42").Value = 81573
$ws.Range("D42").Value = 73402
$ws.Range("E42").Value = 7652
$ws.Range("H42").Value = 519

# Row 61: Suiza - refreshed case counts
$ws.Range("B61").Value = 40645
$ws.Range("C61").Value = 383
$ws.Range("E61").Value = 3842
$ws.Range("G61").Value = 1
$ws.Range("H61").Value = 2003

# Row 86: Senegal - refreshed case counts
$ws.Range("B86").Value = 13186
$ws.Range("C86").Value = 130
$ws.Range("D86").Value = 8852
$ws.Range("E86").Value = 4059
$ws.Range("G86").Value = 1
$ws.Range("H86").Value = 275

# Row 88: country label changes from Zambia to Libia (re-sort), with refreshed case counts
$ws.Range("A88").Value = "Libia"
$ws.Range("B88").Value = 11834
$ws.Range("C88").Value = 553
$ws.Range("D88").Value = 1152
$ws.Range("E88").Value = 10472
$ws.Range("G88").Value = 7
$ws.Range("H88").Value = 210

# Row 89: country label changes from Libia to Zambia (re-sort), with refreshed case counts
$ws.Range("A89").Value = "Zambia"
$ws.Range("B89").Value = 11285
$ws.Range("D89").Value = 10400
$ws.Range("E89").Value = 603
$ws.Range("H89").Value = 282

# Row 101: Finlandia - refreshed case counts
$ws.Range("D101").Value = 7200
$ws.Range("E101").Value = 467

# Row 111: Hong Kong - refreshed case counts
$ws.Range("B111").Value = 4736
$ws.Range("C111").Value = 25
$ws.Range("D111").Value = 4161
$ws.Range("E111").Value = 497

# Row 143: country label changes from Jamaica to Malta (re-sort), with refreshed case counts
$ws.Range("A143").Value = "Malta"
$ws.Range("B143").Value = 1751
$ws.Range("C143").Value = 46
$ws.Range("D143").Value = 1077
$ws.Range("E143").Value = 664
$ws.Range("G143").Value = 0
$ws.Range("H143").Value = 10

# Row 144: country label changes from Jordania to Jamaica (re-sort), with refreshed case counts
$ws.Range("A144").Value = "Jamaica"
$ws.Range("B144").Value = 1732
$ws.Range("C144").Value = 120
$ws.Range("D144").Value = 840
$ws.Range("E144").Value = 873
$ws.Range("G144").Value = 3
$ws.Range("H144").Value = 19

# Row 145: country label changes from Malta to Jordania (re-sort), with refreshed case counts
$ws.Range("A145").Value = "Jordania"
$ws.Range("B145").Value = 1716
$ws.Range("D145").Value = 1344
$ws.Range("E145").Value = 358
$ws.Range("H145").Value = 14

# Row 185: Gibraltar - refreshed case counts
$ws.Range("B185").Value = 270
$ws.Range("C185").Value = 14
$ws.Range("E185").Value = 67
